# Risimulate data using pert distribution
# Update the simulated values in rows 2-7 (columns B..I) to reflect a fresh
# draw from the PERT distribution used to generate this sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-1, 23, 33, 27, 0.4, 0.8, 0.6, 100),
    @(-5, 17, 28, 22, 0.6, 1,   0.7, 81),
    @(0,  16, 22, 21, 0.6, 1.1, 0.8, 60),
    @(-6, 1,  9,  5,  0.6, 1,   0.7, 79),
    @(-2, 5,  9,  8,  0.6, 0.9, 0.8, 68),
    @(6,  11, 20, 15, 0.8, 1.1, 0.9, 67)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
